# Trip_Planner_Pipeline.pptx edit:
#   - Remove the "Where are you currently located?" question (TextBox 2)
#     and its supporting highlight box (Rectangle 7) + arrow connector
#     (Straight Arrow Connector 10).
#   - Reflow "TextBox 5" (the "THIS PAGE TO BE USED ONLY IF DESIRED ..."
#     banner) to a wider box, and reword its text to mention DESTINATION,
#     split across several runs.
#   - Nudge "Rectangle 6" (the page-title highlight) down slightly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- Delete "TextBox 2" ("Where are you currently located?") ---
$s.Shapes.Item("TextBox 2").Delete()

# --- Delete "Rectangle 7" (the highlight box around the removed question) ---
$s.Shapes.Item("Rectangle 7").Delete()

# --- Delete "Straight Arrow Connector 10" (the arrow feeding the removed question) ---
$s.Shapes.Item("Straight Arrow Connector 10").Delete()

# --- Reposition / resize "TextBox 5" and reword its text ---
$tb5 = $s.Shapes.Item("TextBox 5")
$tb5.Left = 100.48181102362204
$tb5.Top = 22.983700787401574
$tb5.Width = 523.560157480315
$tb5.Height = 50.892208099365234

$tr = $tb5.TextFrame.TextRange
$para1 = $tr.Paragraphs(1)
$para1.Text = "THIS PAGE TO BE USED ONLY IF DESIRED "
[void]$para1.InsertAfter("DESTINATION")
[void]$para1.InsertAfter(" ")
[void]$para1.InsertAfter("IS UNKNOWN")

# --- Nudge "Rectangle 6" down slightly (Top only; Left/Width/Height unchanged) ---
$r6 = $s.Shapes.Item("Rectangle 6")
$r6.Top = 24.206615447998047
